$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows 112-121 hold a new batch of "IMDB reviews" GBM/XGBoost/LightGBM/
# CatBoost accuracy results. Column A ("GBM") values are refreshed with the
# latest experiment run; the other columns are unchanged.
$ws.Range("A112").Value = 0.79600000000000004
$ws.Range("A113").Value = 0.83
$ws.Range("A114").Value = 0.83099999999999996
$ws.Range("A115").Value = 0.81499999999999995
$ws.Range("A116").Value = 0.83
$ws.Range("A117").Value = 0.82899999999999996
$ws.Range("A118").Value = 0.79300000000000004
$ws.Range("A119").Value = 0.82399999999999995
$ws.Range("A120").Value = 0.82299999999999995
$ws.Range("A121").Value = 0.82099999999999995

# These rows previously carried a leftover custom style (a bordered/centered
# xf that duplicated the workbook's default "Normalny" style). Revert the
# block back to the default style, same as the rest of the data rows.
$ws.Range("A112:E121").Style = "Normalny"
